$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.265.16"
$ws.Cells.Item(2, 5).Value = "  -4.74%  "

$ws.Cells.Item(3, 4).Value = "2.235.14"
$ws.Cells.Item(3, 5).Value = "  -5.62%  "

$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "315.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.91%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "101.83"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -5.77%  "

$ws.Cells.Item(7, 5).Value = "  -6.69%  "

$ws.Cells.Item(8, 5).Value = "  -0.02%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.561"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -8.01%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "37.06"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -9.16%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.53"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.74%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0825"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -10.01%  "

$ws.Cells.Item(13, 5).Value = "  -9.25%  "

$ws.Cells.Item(14, 5).Value = "  -1.17%  "

$ws.Cells.Item(15, 4).Value = "2.576.62"
$ws.Cells.Item(15, 5).Value = "  -5.56%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.861"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -11.97%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.36"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -6.22%  "

$ws.Cells.Item(18, 4).Value = "2.239.21"
$ws.Cells.Item(18, 5).Value = "  -5.31%  "

$ws.Cells.Item(19, 4).Value = "43.186.98"
$ws.Cells.Item(19, 5).Value = "  -4.96%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.33"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.76%  "

$ws.Cells.Item(21, 5).Value = "  -9.61%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.49"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -10.03%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "65.61"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -10.44%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "3.17"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -10.39%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "238.84"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -7.98%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.12"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -10.95%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.10%  "

$ws.Cells.Item(28, 5).Value = "  +2.24%  "

$ws.Cells.Item(29, 5).Value = "  -9.86%  "

$ws.Cells.Item(30, 5).Value = "  -2.56%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.38"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -12.21%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "20.53"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -7.87%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "34.61"
$ws.Cells.Item(33, 4).Style = "Normal"

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0873"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -10.50%  "

$ws.Cells.Item(35, 5).Value = "  -7.08%  "

$ws.Cells.Item(36, 5).Value = "  -6.08%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.19"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +8.65%  "

$ws.Cells.Item(38, 5).Value = "  +1.91%  "

$ws.Cells.Item(39, 5).Value = "  -6.35%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "4.45"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -5.37%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.104"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -11.26%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.72"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -5.57%  "

$ws.Cells.Item(43, 5).Value = "  -8.53%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "12.90"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.63%  "

$ws.Cells.Item(45, 5).Value = "  -0.03%  "

$ws.Cells.Item(46, 4).Value = "1.800.52"
$ws.Cells.Item(46, 5).Value = "  -0.65%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "89.14"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -9.93%  "

$ws.Cells.Item(48, 5).Value = "  -9.28%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "77.07"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -7.84%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.33"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -8.10%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "59.22"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -15.29%  "
